$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.603.81"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.263.64"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.32%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "3.257.46"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("E10").Value = "  -5.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.566"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.67"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.35%  "
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "691.91"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "3.793.21"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.19"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.36%  "
$ws.Range("D17").Value = "66.701.94"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "3.256.81"
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.11"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.877"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.79"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.45%  "
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.56"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("E26").Value = "  -3.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.64"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.27"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("E29").Value = "  -4.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.26"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.60"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "572.79"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.95%  "
$ws.Range("D33").Value = "3.826.68"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.72"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.45%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.30"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -15.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "54.94"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.36"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("E41").Value = "  -4.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "31.15"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("E43").Value = "  -6.58%  "
$ws.Range("E44").Value = "  -4.38%  "
$ws.Range("E45").Value = "  -7.65%  "
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("E50").Value = "  +3.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.23"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.11%  "
